# Update the price-list column (D) for a couple of client rows and
# move the active selection, per the authored change:
#  - Row 3 (FABIAN FERRETERIA):  LISTA_PRECIOS  D -> E
#  - Row 4 (MARIANO):            LISTA_PRECIOS  D -> E
#  - Row 5 (FORRAJERIA MARTIN):  LISTA_PRECIOS  D -> F
#  - Active cell selection moves from C10 to D6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "E"
$ws.Range("D4").Value = "E"
$ws.Range("D5").Value = "F"

$ws.Range("D6").Select()
